$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old data block (rows 2-23) entirely so new rows are written
# without inheriting the previous ht="17" row-height override.
$ws.Rows("2:23").Delete()

# New quarters added at the top of the table.
$ws.Range("B2").Value = "3T19"
$ws.Range("C2").Value = 414.8
$ws.Range("B3").Value = "4T19"
$ws.Range("C3").Value = 265.1

$ws.Range("B4").Value = "1T20"
$ws.Range("C4").Value = 289.60000000000002
$ws.Range("B5").Value = "2T20"
$ws.Range("C5").Value = 3
$ws.Range("B6").Value = "3T20"
$ws.Range("C6").Value = 62
$ws.Range("B7").Value = "4T20"
$ws.Range("C7").Value = 162.80000000000001
$ws.Range("B8").Value = "1T21"
$ws.Range("C8").Value = 165.9
$ws.Range("B9").Value = "2T21"
$ws.Range("C9").Value = 115.5
$ws.Range("B10").Value = "3T21"
$ws.Range("C10").Value = 230.4
$ws.Range("B11").Value = "4T21"
$ws.Range("C11").Value = 314
$ws.Range("B12").Value = "1T22"
$ws.Range("C12").Value = 292.8
$ws.Range("B13").Value = "2T22"
$ws.Range("C13").Value = 269.7
$ws.Range("B14").Value = "3T22"
$ws.Range("C14").Value = 337.6
$ws.Range("B15").Value = "4T22"
$ws.Range("C15").Value = 321.39999999999998
$ws.Range("B16").Value = "1T23"
$ws.Range("C16").Value = 295.5
$ws.Range("B17").Value = "2T23"
$ws.Range("C17").Value = 269.35700000000003
$ws.Range("B18").Value = "3T23"
$ws.Range("C18").Value = 375.83824480356247
$ws.Range("B19").Value = "4T23"
$ws.Range("C19").Value = 352.24628715989599
$ws.Range("B20").Value = "1T24"
$ws.Range("C20").Value = 317.35410196326967
$ws.Range("B21").Value = "2T24"
$ws.Range("C21").Value = 294.02237340953673
$ws.Range("B22").Value = "3T24"
$ws.Range("C22").Value = 363.81348157635279
$ws.Range("B23").Value = "4T24"
$ws.Range("C23").Value = 366.41134456999987
$ws.Range("B24").Value = "1T25"
$ws.Range("C24").Value = 362.2
$ws.Range("B25").Value = "2T25"
$ws.Range("C25").Value = 341.80872406186722

# ---- Formatting ---------------------------------------------------
$all = $ws.Range("B2:C25")
$all.Font.Name = "Arial"
$all.Font.Bold = $true
$all.Font.Size = 10

$all.Interior.Pattern = 1
$all.Interior.ThemeColor = 2
$all.Interior.TintAndShade = 0

$all.HorizontalAlignment = -4108
$all.VerticalAlignment = -4108

$ws.Range("B4:B25").WrapText = $true

$ws.Range("G13").Select()
